$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" updates ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 336.4959975
$schedule.Range("F2").Value = 7.418342096560846
$schedule.Range("E3").Value = 446.4720975
$schedule.Range("F3").Value = 29.52857787698413

# --- Sheet "Detailed" updates ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B13").Value = 65
$detailed.Range("B14").Value = 76.36297999999999
$detailed.Range("B15").Value = 77.94
$detailed.Range("C15").Value = "historical"
$detailed.Range("B16").Value = 40.54
$detailed.Range("C16").Value = "historical"
$detailed.Range("B18").Value = -1.06158
$detailed.Range("B19").Value = 0.00967
$detailed.Range("B20").Value = -5.01
$detailed.Range("B22").Value = -6.66927
$detailed.Range("B23").Value = -6.83948
$detailed.Range("B24").Value = -7.94171
$detailed.Range("B25").Value = -8.161899999999999
$detailed.Range("B26").Value = -6.5522
$detailed.Range("B27").Value = -7.86159
$detailed.Range("B28").Value = -8.008850000000001
$detailed.Range("B29").Value = -7.43552
$detailed.Range("B30").Value = -7.21509
$detailed.Range("B31").Value = -7.85784
$detailed.Range("B32").Value = -6.42637
$detailed.Range("B33").Value = -6.21778
$detailed.Range("B34").Value = -3.0719
$detailed.Range("B35").Value = -5.16056
$detailed.Range("B36").Value = -2.76655
$detailed.Range("B37").Value = 9.643750000000001
$detailed.Range("B38").Value = 9.84629
$detailed.Range("B39").Value = 33.09616
$detailed.Range("B40").Value = 56.98
$detailed.Range("B41").Value = 56.87831
$detailed.Range("B42").Value = 57.3
$detailed.Range("B44").Value = 57.3
$detailed.Range("B45").Value = 57.3
$detailed.Range("B46").Value = 57.06007
